$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("s69","s69_e29_24019-32_1_1.jpeg","meltpatch","1047","1581","10","10","9","1"),
  @("s70","s70_e30_24019-32_1_2.jpeg","meltpatch","1407","370","21","10","23","2"),
  @("s71","s71_e31_24019-32_1_3.jpeg","meltpatch","1946","514","31","10","102","3"),
  @("s72","s72_e32_24019-32_1_4.jpeg","meltpatch","398","1095","42","10","12","4"),
  @("s73","s73_e33_24019-32_1_5.jpeg","meltpatch","1580","921","52","10","39","5"),
  @("s74","s74_e34_24019-32_1_6.jpeg","meltpatch","1546","1494","63","10","47","6"),
  @("s75","s75_e35_24019-32_1_7.jpeg","meltpatch","1839","37","73","10","72","7"),
  @("s76","s76_e36_24019-32_1_8.jpeg","meltpatch","839","789","84","10","139","8"),
  @("s77","s77_e37_24019-32_2_7.jpeg","meltpatch","1009","859","94","10","99","9"),
  @("s78","s78_e38_24019-32_2_6.jpeg","meltpatch","234","702","105","10","167","10"),
  @("s79","s79_e39_24019-32_2_5.jpeg","meltpatch","847","1585","157","10","85","15"),
  @("s80","s80_e40_24019-32_2_4.jpeg","meltpatch","371","1023","210","10","9","20"),
  @("s81","s81_e41_24019-32_2_3.jpeg","meltpatch","2065","1558","15","15","35","1"),
  @("s82","s82_e42_24019-32_2_2.jpeg","meltpatch","774","204","31","15","103","2"),
  @("s83","s83_e43_24019-32_2_1.jpeg","meltpatch","677","233","47","15","86","3"),
  @("s84","s84_e44_24019-32_2_0.jpeg","meltpatch","692","1152","63","15","21","4"),
  @("s85","s85_e45_24019-32_3_1.jpeg","meltpatch","1794","648","78","15","55","5"),
  @("s86","s86_e46_24019-32_3_2.jpeg","meltpatch","1864","850","94","15","123","6"),
  @("s87","s87_e47_24019-32_3_3.jpeg","meltpatch","2009","682","110","15","153","7"),
  @("s88","s88_e48_24019-32_3_4.jpeg","meltpatch","2096","1058","126","15","116","8"),
  @("s89","s89_e49_24019-32_3_5.jpeg","meltpatch","144","1319","142","15","12","9"),
  @("s90","s90_e50_24019-32_3_6.jpeg","meltpatch","1765","1266","157","15","85","10"),
  @("s91","s91_e51_24019-32_3_7.jpeg","meltpatch","1779","128","205","15","97","13"),
  @("s92","s92_e52_24019-32_3_8.jpeg","meltpatch","474","1557","21","21","87","1"),
  @("s93","s93_e53_24019-32_4_7.jpeg","meltpatch","734","96","42","21","164","2"),
  @("s94","s94_e54_24019-32_4_6.jpeg","meltpatch","60","323","63","21","83","3"),
  @("s95","s95_e55_24019-32_4_5.jpeg","meltpatch","423","402","84","21","157","4"),
  @("s96","s96_e56_24019-32_4_4.jpeg","meltpatch","979","832","105","21","118","5"),
  @("s97","s97_e57_24019-32_4_3.jpeg","meltpatch","1405","366","126","21","55","6"),
  @("s98","s98_e58_24019-32_4_2.jpeg","meltpatch","588","1446","147","21","43","7"),
  @("s99","s99_e59_24019-32_4_1.jpeg","meltpatch","1543","748","168","21","29","8"),
  @("s100","s100_e60_24019-32_4_0.jpeg","meltpatch","1375","36","189","21","143","9"),
  @("s101","s101_e61_24019-32_5_1.jpeg","meltpatch","1701","476","210","21","172","10"),
  @("s102","s102_e62_24019-32_5_2.jpeg","meltpatch","1835","1351","26","26","87","1")
)

$startRow = 70
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $val = $rowData[$j]
        if ($col -ge 4) {
            $ws.Cells.Item($r, $col).Value = "'" + $val
        } else {
            $ws.Cells.Item($r, $col).Value = $val
        }
    }
}